$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137 (sheet ALC)
$ws.Range("H137").Value = 1685.6957
$ws.Range("I137").Value = 1671.409
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 5014.227000000001
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -2464.227000000001
$ws.Range("N137").Value = -11100


$ws = $wb.Worksheets.Item("ARM")
# Row 122 (sheet ARM)
$ws.Range("H122").Value = 2185.1428
$ws.Range("I122").Value = 1888.6666
$ws.Range("J122").Value = 2718.8
$ws.Range("K122").Value = 5665.9998
$ws.Range("L122").Value = 8156.400000000001
$ws.Range("M122").Value = -3215.9998
$ws.Range("N122").Value = -13056.4

# Row 132 (sheet ARM)
$ws.Range("H132").Value = 27313.625
$ws.Range("I132").Value = 21702.2
$ws.Range("J132").Value = 36666
$ws.Range("K132").Value = 65106.60000000001
$ws.Range("L132").Value = 109998
$ws.Range("M132").Value = -62576.60000000001
$ws.Range("N132").Value = -115058


$ws = $wb.Worksheets.Item("BSM")
# Row 20 (sheet BSM)
$ws.Range("H20").Value = 2048.7334
$ws.Range("I20").Value = 1818
$ws.Range("J20").Value = 2683.25
$ws.Range("K20").Value = 1818
$ws.Range("L20").Value = 2683.25
$ws.Range("M20").Value = -1571
$ws.Range("N20").Value = -3177.25

# Row 134 (sheet BSM)
$ws.Range("H134").Value = 10281.454
$ws.Range("I134").Value = 11239.6
$ws.Range("J134").Value = 700
$ws.Range("K134").Value = 33718.8
$ws.Range("L134").Value = 2100
$ws.Range("M134").Value = -31183.8
$ws.Range("N134").Value = -7170


$ws = $wb.Worksheets.Item("CRP")
# Row 94 (sheet CRP)
$ws.Range("H94").Value = 402138.9
$ws.Range("I94").Value = 672104
$ws.Range("J94").Value = 286439.56
$ws.Range("K94").Value = 672104
$ws.Range("L94").Value = 286439.56
$ws.Range("M94").Value = -671653
$ws.Range("N94").Value = -287341.56

# Row 99 (sheet CRP)
$ws.Range("H99").Value = 1637.3334
$ws.Range("I99").Value = 1416.875
$ws.Range("J99").Value = 2078.25
$ws.Range("K99").Value = 1416.875
$ws.Range("L99").Value = 2078.25
$ws.Range("M99").Value = 81.125
$ws.Range("N99").Value = -5074.25

# Row 105 (sheet CRP)
$ws.Range("H105").Value = 2866.25
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 3488.3333
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 3488.3333
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -6982.3333

# Row 122 (sheet CRP)
$ws.Range("H122").Value = 2414
$ws.Range("I122").Value = 2464.5715
$ws.Range("J122").Value = 2355
$ws.Range("K122").Value = 7393.7145
$ws.Range("L122").Value = 7065
$ws.Range("M122").Value = -4943.7145
$ws.Range("N122").Value = -11965

# Row 126 (sheet CRP)
$ws.Range("H126").Value = 1637.3334
$ws.Range("I126").Value = 1416.875
$ws.Range("J126").Value = 2078.25
$ws.Range("K126").Value = 4250.625
$ws.Range("L126").Value = 6234.75
$ws.Range("M126").Value = -1780.625
$ws.Range("N126").Value = -11174.75


$ws = $wb.Worksheets.Item("GSM")
# Row 43 (sheet GSM)
$ws.Range("H43").Value = 3696.5
$ws.Range("I43").Value = 3696.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3696.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3545.5
$ws.Range("N43").ClearContents()

# Row 46 (sheet GSM)
$ws.Range("H46").Value = 26965
$ws.Range("I46").Value = 10525
$ws.Range("J46").Value = 29705
$ws.Range("K46").Value = 10525
$ws.Range("L46").Value = 29705
$ws.Range("M46").Value = -10369
$ws.Range("N46").Value = -30017

# Row 57 (sheet GSM)
$ws.Range("H57").Value = 30061
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 30061
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 30061
$ws.Range("N57").Value = -31701
$ws.Range("M57").ClearContents()

# Row 80 (sheet GSM)
$ws.Range("H80").Value = 2689.2307
$ws.Range("I80").Value = 2611.25
$ws.Range("J80").Value = 2814
$ws.Range("K80").Value = 2611.25
$ws.Range("L80").Value = 2814
$ws.Range("M80").Value = -1613.25
$ws.Range("N80").Value = -4810

# Row 83 (sheet GSM)
$ws.Range("H83").Value = 2689.2307
$ws.Range("I83").Value = 2611.25
$ws.Range("J83").Value = 2814
$ws.Range("K83").Value = 13056.25
$ws.Range("L83").Value = 14070
$ws.Range("M83").Value = -8064.25
$ws.Range("N83").Value = -24054

# Row 102 (sheet GSM)
$ws.Range("H102").Value = 2016.1538
$ws.Range("I102").Value = 2016.1538
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2016.1538
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -394.1538
$ws.Range("N102").ClearContents()

# Row 122 (sheet GSM)
$ws.Range("H122").Value = 1512.1428
$ws.Range("I122").Value = 957
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 2871
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -421
$ws.Range("N122").Value = -13600


$ws = $wb.Worksheets.Item("LTW")
# Row 7 (sheet LTW)
$ws.Range("H7").Value = 3257.7778
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 3386.6667
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3386.6667
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -3610.6667

# Row 126 (sheet LTW)
$ws.Range("H126").Value = 3257.7778
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3386.6667
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 10160.0001
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -15100.0001

# Row 132 (sheet LTW)
$ws.Range("H132").Value = 33684.812
$ws.Range("I132").Value = 40240.04
$ws.Range("J132").Value = 5278.8335
$ws.Range("K132").Value = 120720.12
$ws.Range("L132").Value = 15836.5005
$ws.Range("M132").Value = -118190.12
$ws.Range("N132").Value = -20896.5005

# Row 136 (sheet LTW)
$ws.Range("H136").Value = 2583.7104
$ws.Range("I136").Value = 1790.2609
$ws.Range("J136").Value = 3800.3333
$ws.Range("K136").Value = 5370.7827
$ws.Range("L136").Value = 11400.9999
$ws.Range("M136").Value = -2820.7827
$ws.Range("N136").Value = -16500.9999


$ws = $wb.Worksheets.Item("WVR")
# Row 52 (sheet WVR)
$ws.Range("H52").Value = 27257.75
$ws.Range("I52").Value = 6992
$ws.Range("J52").Value = 47523.5
$ws.Range("K52").Value = 6992
$ws.Range("L52").Value = 47523.5
$ws.Range("M52").Value = -6766
$ws.Range("N52").Value = -47975.5

# Row 132 (sheet WVR)
$ws.Range("H132").Value = 10247.083
$ws.Range("I132").Value = 11296.7
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 33890.10000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -31360.10000000001
$ws.Range("N132").Value = -20057

# Row 136 (sheet WVR)
$ws.Range("H136").Value = 76932730
$ws.Range("I136").Value = 100011550
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 300034650
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -300032100
$ws.Range("N136").Value = -15099.9999

